# "Completed Question 8 no bonus"
#
# The metro_budget sheet's "Question 8" block (rows 88-93) asks: for each of
# three ranks entered in row 89 (columns B, D, F), look up which department
# holds that rank for FY17 / FY18 / FY19 and report its diff_pct.
#
#   Row 91 -> FY17  (rank lookup in F2:F52, diff_pct lookup in E2:E52)
#   Row 92 -> FY18  (rank lookup in K2:K52, diff_pct lookup in J2:J52)
#   Row 93 -> FY19  (rank lookup in P2:P52, diff_pct lookup in O2:O52)
#
# For each row/FY pair: col B/D/F = department name for the requested rank
# (B$89/D$89/F$89) via XLOOKUP against the rank column; col C/E/G = that
# department's diff_pct via a second XLOOKUP against column A (Department).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metro_budget")

# Row 91 - FY17 (rank column F, diff_pct column E)
$ws.Range("B91").Formula = "=_xlfn.XLOOKUP(B`$89,`$F`$2:`$F`$52,`$A`$2:`$A`$52)"
$ws.Range("C91").Formula = "=_xlfn.XLOOKUP(B91,`$A`$2:`$A`$52,`$E`$2:`$E`$52)"
$ws.Range("D91").Formula = "=_xlfn.XLOOKUP(D`$89,`$F`$2:`$F`$52,`$A`$2:`$A`$52)"
$ws.Range("E91").Formula = "=_xlfn.XLOOKUP(D91,`$A`$2:`$A`$52,`$E`$2:`$E`$52)"
$ws.Range("F91").Formula = "=_xlfn.XLOOKUP(F`$89,`$F`$2:`$F`$52,`$A`$2:`$A`$52)"
$ws.Range("G91").Formula = "=_xlfn.XLOOKUP(F91,`$A`$2:`$A`$52,`$E`$2:`$E`$52)"

# Row 92 - FY18 (rank column K, diff_pct column J)
$ws.Range("B92").Formula = "=_xlfn.XLOOKUP(B`$89,`$K`$2:`$K`$52,`$A`$2:`$A`$52)"
$ws.Range("C92").Formula = "=_xlfn.XLOOKUP(B92,`$A`$2:`$A`$52,`$J`$2:`$J`$52)"
$ws.Range("D92").Formula = "=_xlfn.XLOOKUP(D`$89,`$K`$2:`$K`$52,`$A`$2:`$A`$52)"
$ws.Range("E92").Formula = "=_xlfn.XLOOKUP(D92,`$A`$2:`$A`$52,`$J`$2:`$J`$52)"
$ws.Range("F92").Formula = "=_xlfn.XLOOKUP(F`$89,`$K`$2:`$K`$52,`$A`$2:`$A`$52)"
$ws.Range("G92").Formula = "=_xlfn.XLOOKUP(F92,`$A`$2:`$A`$52,`$J`$2:`$J`$52)"

# Row 93 - FY19 (rank column P, diff_pct column O)
$ws.Range("B93").Formula = "=_xlfn.XLOOKUP(B`$89,`$P`$2:`$P`$52,`$A`$2:`$A`$52)"
$ws.Range("C93").Formula = "=_xlfn.XLOOKUP(B93,`$A`$2:`$A`$52,`$O`$2:`$O`$52)"
$ws.Range("D93").Formula = "=_xlfn.XLOOKUP(D`$89,`$P`$2:`$P`$52,`$A`$2:`$A`$52)"
$ws.Range("E93").Formula = "=_xlfn.XLOOKUP(D93,`$A`$2:`$A`$52,`$O`$2:`$O`$52)"
$ws.Range("F93").Formula = "=_xlfn.XLOOKUP(F`$89,`$P`$2:`$P`$52,`$A`$2:`$A`$52)"
$ws.Range("G93").Formula = "=_xlfn.XLOOKUP(F93,`$A`$2:`$A`$52,`$O`$2:`$O`$52)"

# The diff_pct cells already carried the percentage style before this edit
# (style index 5 / "0.00%"); re-assert it explicitly so it's preserved no
# matter how the formulas were written in.
$ws.Range("C91,E91,G91,C92,E92,G92,C93,E93,G93").NumberFormat = "0.00%"

# Reflect where the author ended up after finishing the exercise: scrolled
# down so row 76 is at the top of the window, with B98 the active selection.
$ws.Application.Goto($ws.Range("A76"))
$excel.ActiveWindow.ScrollRow = 76
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B98").Select()
